# Update "total views/likes" figures in column F across the four sheets,
# mirroring the content refresh recorded at commit 456a3b4 for the
# gh-pages generated output.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 881
$ws.Cells.Item(3, 6).Value = 1007
$ws.Cells.Item(4, 6).Value = 786
$ws.Cells.Item(5, 6).Value = 866
$ws.Cells.Item(6, 6).Value = 444
$ws.Cells.Item(7, 6).Value = 684
$ws.Cells.Item(8, 6).Value = 158
$ws.Cells.Item(9, 6).Value = 1282
$ws.Cells.Item(10, 6).Value = 712
$ws.Cells.Item(11, 6).Value = 413
$ws.Cells.Item(12, 6).Value = 545
$ws.Cells.Item(15, 6).Value = 960
$ws.Cells.Item(16, 6).Value = 16
$ws.Cells.Item(17, 6).Value = 402
$ws.Cells.Item(20, 6).Value = 584
$ws.Cells.Item(21, 6).Value = 146
$ws.Cells.Item(22, 6).Value = 633
$ws.Cells.Item(24, 6).Value = 991

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 339
$ws.Cells.Item(3, 6).Value = 111
$ws.Cells.Item(7, 6).Value = 241
$ws.Cells.Item(8, 6).Value = 55
$ws.Cells.Item(11, 6).Value = 112

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 376

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 376
$ws.Cells.Item(3, 6).Value = 339
$ws.Cells.Item(4, 6).Value = 881
$ws.Cells.Item(5, 6).Value = 1007
$ws.Cells.Item(6, 6).Value = 786
$ws.Cells.Item(7, 6).Value = 866
$ws.Cells.Item(8, 6).Value = 444
$ws.Cells.Item(9, 6).Value = 684
$ws.Cells.Item(10, 6).Value = 158
$ws.Cells.Item(11, 6).Value = 1282
$ws.Cells.Item(12, 6).Value = 712
$ws.Cells.Item(13, 6).Value = 111
$ws.Cells.Item(15, 6).Value = 413
$ws.Cells.Item(16, 6).Value = 545
$ws.Cells.Item(20, 6).Value = 960
$ws.Cells.Item(22, 6).Value = 16
$ws.Cells.Item(23, 6).Value = 402
$ws.Cells.Item(26, 6).Value = 241
$ws.Cells.Item(27, 6).Value = 55
$ws.Cells.Item(28, 6).Value = 584
$ws.Cells.Item(31, 6).Value = 112
$ws.Cells.Item(32, 6).Value = 112
$ws.Cells.Item(33, 6).Value = 146
$ws.Cells.Item(34, 6).Value = 633
$ws.Cells.Item(36, 6).Value = 991
